# tested and refactored job application method
# Rewrites the single-row job-application log into a multi-row history
# (6 rows of data instead of 2), reusing the same account (email/password)
# against five different job postings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks before rewriting cell values/layout - they will
# be re-created (in a new order) once the final data is in place.
$ws.Cells.Hyperlinks.Delete()

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value2 = "Email"
$ws.Range("B1").Value2 = "Parola"
$ws.Range("C1").Value2 = "Job"
$ws.Range("D1").Value2 = "Domeniu"
$ws.Range("E1").Value2 = "Oras"
$ws.Range("F1").Value2 = "Experience"
$ws.Range("G1").Value2 = "Job Type"

# --- Data rows ----------------------------------------------------------
# Row 2: secretar / Secretariat-Administrativ / Bucuresti
$ws.Range("A2").Value2 = "beatrice.dobre@asmi.ro"
$ws.Range("B2").Value2 = "Parola2000!"
$ws.Range("C2").Value2 = "secretar"
$ws.Range("D2").Value2 = "Secretariat-Administrativ"
$ws.Range("E2").Value2 = "Bucuresti"
$ws.Range("F2").Value2 = "0 - 1 an experienta,1 - 5 ani experienta"
$ws.Range("G2").Value2 = "full-time,practica"

# Row 3: legal / Juridic / bucuresti
$ws.Range("A3").Value2 = "beatrice.dobre@asmi.ro"
$ws.Range("B3").Value2 = "Parola2000!"
$ws.Range("C3").Value2 = "legal"
$ws.Range("D3").Value2 = "Juridic"
$ws.Range("E3").Value2 = "bucuresti"
$ws.Range("F3").Value2 = "1 - 5 ani experienta"
$ws.Range("G3").Value2 = "full-time"

# Row 4: inginer / Inginerie / Timisoara
$ws.Range("A4").Value2 = "beatrice.dobre@asmi.ro"
$ws.Range("B4").Value2 = "Parola2000!"
$ws.Range("C4").Value2 = "inginer"
$ws.Range("D4").Value2 = "Inginerie"
$ws.Range("E4").Value2 = "Timisoara"
$ws.Range("F4").Value2 = "peste 5 ani experienta,Manager"
$ws.Range("G4").Value2 = "full-time"

# Row 5: profesor / Educatie-Training / Bucuresti
$ws.Range("A5").Value2 = "beatrice.dobre@asmi.ro"
$ws.Range("B5").Value2 = "Parola2000!"
$ws.Range("C5").Value2 = "profesor"
$ws.Range("D5").Value2 = "Educatie-Training"
$ws.Range("E5").Value2 = "Bucuresti"
$ws.Range("F5").Value2 = "Student--Absolvent"
$ws.Range("G5").Value2 = "part-time"

# Row 6: IT / IT-Software / bucuresti (the original application row)
$ws.Range("A6").Value2 = "beatrice.dobre@asmi.ro"
$ws.Range("B6").Value2 = "Parola2000!"
$ws.Range("C6").Value2 = "IT"
$ws.Range("D6").Value2 = "IT-Software"
$ws.Range("E6").Value2 = "bucuresti"
$ws.Range("F6").Value2 = "0 - 1 an experienta"
$ws.Range("G6").Value2 = "full-time"

# --- Re-create the mailto hyperlinks (+ styling) for column A, rows 2-6 -
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:beatrice.dobre@asmi.ro") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:beatrice.dobre@asmi.ro") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:beatrice.dobre@asmi.ro") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:beatrice.dobre@asmi.ro") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:beatrice.dobre@asmi.ro") | Out-Null

$ws.Range("A2:A6").Style = "Hyperlink"

# --- Selection, matching the saved view state ---------------------------
$ws.Range("D14").Select() | Out-Null
